$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" email lists (text-only changes) ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- Update summary statistics numbers ---
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 20
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 20

# --- Row 10 (HISTOLOGY session 2) moves from "Pending" to "Not Recorded" ---
$ws.Range("I10").Value = "Not Recorded"

# Copy the visual format of row 29 (already styled as "Not Recorded") onto row 10
# so the same underlying style entry (font/fill) is reused instead of creating a
# brand-new one.
$ws.Range("A29:I29").Copy() | Out-Null
$ws.Range("A10:I10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
